$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1049.28
$ws.Range("J112").Value = 1072.1666
$ws.Range("L112").Value = 3216.4998
$ws.Range("N112").Value = -5432.4998
$ws.Range("H132").Value = 4550742.5
$ws.Range("I132").Value = 5005566.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 15016699.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -15014169.5
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 58800
$ws.Range("J136").Value = 58800
$ws.Range("L136").Value = 58800
$ws.Range("N136").Value = -69000
$ws.Range("H138").Value = 2701.203
$ws.Range("I138").Value = 1621.75
$ws.Range("J138").Value = 3141.796
$ws.Range("K138").Value = 4865.25
$ws.Range("L138").Value = 9425.387999999999
$ws.Range("M138").Value = 274.75
$ws.Range("N138").Value = -19705.388
$ws.Range("H139").Value = 69774.5
$ws.Range("J139").Value = 69774.5
$ws.Range("L139").Value = 69774.5
$ws.Range("N139").Value = -80054.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35550.38
$ws.Range("I2").Value = 1132.875
$ws.Range("J2").Value = 200754.4
$ws.Range("K2").Value = 1132.875
$ws.Range("L2").Value = 200754.4
$ws.Range("M2").Value = -1019.875
$ws.Range("N2").Value = -200980.4
$ws.Range("H116").Value = 35550.38
$ws.Range("I116").Value = 1132.875
$ws.Range("J116").Value = 200754.4
$ws.Range("K116").Value = 1132.875
$ws.Range("L116").Value = 200754.4
$ws.Range("M116").Value = 1161.125
$ws.Range("N116").Value = -205342.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35550.38
$ws.Range("I3").Value = 1132.875
$ws.Range("J3").Value = 200754.4
$ws.Range("K3").Value = 1132.875
$ws.Range("L3").Value = 200754.4
$ws.Range("M3").Value = -1018.875
$ws.Range("N3").Value = -200982.4
$ws.Range("H86").Value = 37176.03
$ws.Range("I86").Value = 45593.12
$ws.Range("J86").Value = 2104.8333
$ws.Range("K86").Value = 45593.12
$ws.Range("L86").Value = 2104.8333
$ws.Range("M86").Value = -44470.12
$ws.Range("N86").Value = -4350.8333
$ws.Range("H89").Value = 37176.03
$ws.Range("I89").Value = 45593.12
$ws.Range("J89").Value = 2104.8333
$ws.Range("K89").Value = 227965.6
$ws.Range("L89").Value = 10524.1665
$ws.Range("M89").Value = -222349.6
$ws.Range("N89").Value = -21756.1665
$ws.Range("H134").Value = 12008.167
$ws.Range("I134").Value = 13388.366
$ws.Range("J134").Value = 3924.1428
$ws.Range("K134").Value = 40165.098
$ws.Range("L134").Value = 11772.4284
$ws.Range("M134").Value = -37630.098
$ws.Range("N134").Value = -16842.4284
$ws.Range("H138").Value = 70860
$ws.Range("J138").Value = 70860
$ws.Range("L138").Value = 70860
$ws.Range("N138").Value = -81140
$ws.Range("H140").Value = 48944.445
$ws.Range("J140").Value = 48944.445
$ws.Range("L140").Value = 48944.445
$ws.Range("N140").Value = -59304.445
$ws.Range("H141").Value = 59500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 59500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 59500
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -69860

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 746.86664
$ws.Range("J16").Value = 2450
$ws.Range("L16").Value = 2450
$ws.Range("N16").Value = -3024
$ws.Range("H113").Value = 746.86664
$ws.Range("J113").Value = 2450
$ws.Range("L113").Value = 2450
$ws.Range("N113").Value = -6790
$ws.Range("H138").Value = 72882.22
$ws.Range("J138").Value = 72882.22
$ws.Range("L138").Value = 72882.22
$ws.Range("N138").Value = -83162.22
$ws.Range("H140").Value = 59590
$ws.Range("J140").Value = 59590
$ws.Range("L140").Value = 59590
$ws.Range("N140").Value = -69950

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1383.6666
$ws.Range("I68").Value = 721.5
$ws.Range("J68").Value = 2045.8334
$ws.Range("K68").Value = 2164.5
$ws.Range("L68").Value = 6137.5002
$ws.Range("M68").Value = -1353.5
$ws.Range("N68").Value = -7759.5002
$ws.Range("H71").Value = 1383.6666
$ws.Range("I71").Value = 721.5
$ws.Range("J71").Value = 2045.8334
$ws.Range("K71").Value = 6493.5
$ws.Range("L71").Value = 18412.5006
$ws.Range("M71").Value = -2437.5
$ws.Range("N71").Value = -26524.5006
$ws.Range("H107").Value = 316347.53
$ws.Range("I107").Value = 481.46155
$ws.Range("J107").Value = 448807.47
$ws.Range("K107").Value = 1444.38465
$ws.Range("L107").Value = 1346422.41
$ws.Range("M107").Value = 475.61535
$ws.Range("N107").Value = -1350262.41

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1749
$ws.Range("I122").Value = 1638.1538
$ws.Range("J122").Value = 1954.8572
$ws.Range("K122").Value = 4914.4614
$ws.Range("L122").Value = 5864.571599999999
$ws.Range("M122").Value = -2464.4614
$ws.Range("N122").Value = -10764.5716
$ws.Range("H138").Value = 78373.75
$ws.Range("J138").Value = 76998.57000000001
$ws.Range("L138").Value = 76998.57000000001
$ws.Range("N138").Value = -87278.57000000001
$ws.Range("H140").Value = 104060
$ws.Range("J140").Value = 104060
$ws.Range("L140").Value = 104060
$ws.Range("N140").Value = -114420

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2164.5293
$ws.Range("I82").Value = 1622
$ws.Range("J82").Value = 2460.4546
$ws.Range("K82").Value = 1622
$ws.Range("L82").Value = 2460.4546
$ws.Range("M82").Value = -1261
$ws.Range("N82").Value = -3182.4546
$ws.Range("H85").Value = 2164.5293
$ws.Range("I85").Value = 1622
$ws.Range("J85").Value = 2460.4546
$ws.Range("K85").Value = 1622
$ws.Range("L85").Value = 2460.4546
$ws.Range("M85").Value = -374
$ws.Range("N85").Value = -4956.4546
$ws.Range("H132").Value = 5478.5
$ws.Range("I132").Value = 5787.7
$ws.Range("J132").Value = 5169.3
$ws.Range("K132").Value = 17363.1
$ws.Range("L132").Value = 15507.9
$ws.Range("M132").Value = -14833.1
$ws.Range("N132").Value = -20567.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 521.1
$ws.Range("I113").Value = 244.42857
$ws.Range("J113").Value = 1166.6666
$ws.Range("K113").Value = 733.28571
$ws.Range("L113").Value = 3499.9998
$ws.Range("M113").Value = 1436.71429
$ws.Range("N113").Value = -7839.9998
$ws.Range("H122").Value = 2135.1
$ws.Range("I122").Value = 1615.591
$ws.Range("J122").Value = 3563.75
$ws.Range("K122").Value = 4846.772999999999
$ws.Range("L122").Value = 10691.25
$ws.Range("M122").Value = -2396.772999999999
$ws.Range("N122").Value = -15591.25
$ws.Range("H138").Value = 45930
$ws.Range("J138").Value = 45930
$ws.Range("L138").Value = 45930
$ws.Range("N138").Value = -56210
$ws.Range("H139").Value = 65114
$ws.Range("J139").Value = 65114
$ws.Range("L139").Value = 65114
$ws.Range("N139").Value = -75394
